$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290, shifting existing rows 290-390 down to 291-391.
$ws.Rows(290).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(290, 1).Value = 4
$ws.Cells.Item(290, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(290, 3).Value = "Los Lagos"
$ws.Cells.Item(290, 4).Value = 44900
$ws.Cells.Item(290, 5).Value = 10
$ws.Cells.Item(290, 6).Value = 100112045
$ws.Cells.Item(290, 7).Value = "Zapallo"
$ws.Cells.Item(290, 8).Value = "Paine"
$ws.Cells.Item(290, 9).Value = "1a nueva(o)"
$ws.Cells.Item(290, 10).Value = 500
$ws.Cells.Item(290, 11).Value = 900
$ws.Cells.Item(290, 12).Value = 900
$ws.Cells.Item(290, 13).Value = 900
$ws.Cells.Item(290, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(290, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(290, 16).Value = 900
$ws.Cells.Item(290, 17).Value = 1
$ws.Cells.Item(290, 18).Value = "Hortaliza"
